$wb = $excel.ActiveWorkbook

# --- Metadata sheet content updates -----------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 0.1.15-beta -> 0.1.16-beta
$ws.Range("B3").Value = "0.1.16-beta"

# Date bump: 2023-06-07T11:47:17-05:00 -> 2023-06-13T11:38:47-05:00
$ws.Range("B8").Value = "2023-06-13T11:38:47-05:00"

# --- Formatting fix -----------------------------------------------------
# The header/body cell styles already carried an <alignment vertical="top"
# wrapText="true"/> child but were missing applyAlignment="true", so Excel
# never actually applied the wrap/vertical-top alignment. Re-apply wrap
# text + top vertical alignment on every sheet's used range so the
# alignment is genuinely flagged as applied (applyAlignment="true").
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $used.WrapText = $true
    $used.VerticalAlignment = -4160  # xlVAlignTop
}
